# Apply the "policy approver and corporate management" sheet additions.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the five new reference-data sheets, in order, right after
#    "Sheet3" (so final order becomes:
#    Sheet1, Sheet3, Branch, Department, Designation, EmployeeGrade,
#    Passportvisaalert, Sheet2, Division, username)
# ---------------------------------------------------------------------
$after = $wb.Worksheets.Item("Sheet3")
$newNames = @("Branch", "Department", "Designation", "EmployeeGrade", "Passportvisaalert")
foreach ($n in $newNames) {
    $s = $wb.Worksheets.Add($null, $after)
    $s.Name = $n
    $after = $s
}

# ---------------------------------------------------------------------
# 2. Branch sheet
# ---------------------------------------------------------------------
$branch = $wb.Worksheets.Item("Branch")

$branchHeaders = @("Telephone", "Mobile", "fax", "Contactname", "ContactEmailAdresss", "ContactMobile", "country", "Postcode", "State", "City", "Address", "Branchemail", "branchname", "division")
for ($i = 0; $i -lt $branchHeaders.Count; $i++) {
    $cell = $branch.Cells.Item(1, $i + 1)
    $cell.Value = $branchHeaders[$i]
    $cell.Font.Bold = $true
}

$branch.Cells.Item(2, 1).Value = 9874651222
$branch.Cells.Item(2, 2).Value = 2165165141
$branch.Cells.Item(2, 4).Value = "anurag upadhyay"
$branch.Cells.Item(2, 5).Value = "anurag.upadhyay@quadlabs.com"
$branch.Cells.Item(2, 6).Value = 9879846551
$branch.Cells.Item(2, 7).Value = "India"
$branch.Cells.Item(2, 8).Value = 877545
$branch.Cells.Item(2, 9).Value = "haryana"
$branch.Cells.Item(2, 10).Value = "gurgaon"
$branch.Cells.Item(2, 11).Value = "720 b jmd megapolis "
$branch.Cells.Item(2, 12).Value = "qa@quadlabs.com"
$branch.Cells.Item(2, 13).Value = "India"
$branch.Cells.Item(2, 14).Value = "DELHI"

$branch.Range("A1:N2").Borders.LineStyle = 1

$branch.Hyperlinks.Add($branch.Cells.Item(2, 5), "mailto:anurag.upadhyay@quadlabs.com")
$branch.Hyperlinks.Add($branch.Cells.Item(2, 12), "mailto:qa@quadlabs.com")
$branch.Cells.Item(2, 5).Borders.LineStyle = 1
$branch.Cells.Item(2, 12).Borders.LineStyle = 1

$branch.Columns.Item(1).ColumnWidth = 11.43
$branch.Columns.Item(2).ColumnWidth = 11
$branch.Columns.Item(3).ColumnWidth = 4.29
$branch.Columns.Item(4).ColumnWidth = 13.86
$branch.Columns.Item(5).ColumnWidth = 29.71
$branch.Columns.Item(6).ColumnWidth = 14.71
$branch.Columns.Item(12).ColumnWidth = 12
$branch.Columns.Item(13).ColumnWidth = 11.43

$branch.PageSetup.Orientation = 1
$branch.Range("G2").Select()

# ---------------------------------------------------------------------
# 3. Department sheet
# ---------------------------------------------------------------------
$dept = $wb.Worksheets.Item("Department")

$deptHeaders = @("Branch", "Departmentname", "DepartmentEmail", "DepartmentPhone", "DepartmentFax", "TravelBudget", "Noofstaff", "Nooftraveller")
for ($i = 0; $i -lt $deptHeaders.Count; $i++) {
    $cell = $dept.Cells.Item(1, $i + 1)
    $cell.Value = $deptHeaders[$i]
    $cell.Font.Bold = $true
}

$dept.Cells.Item(2, 1).Value = "India"
$dept.Cells.Item(2, 2).Value = "finance"
$dept.Cells.Item(2, 3).Value = "finance@bmw.in"
$dept.Cells.Item(2, 4).Value = 9999999999
$dept.Cells.Item(2, 5).Value = 9999999999
$dept.Cells.Item(2, 6).Value = 500000
$dept.Cells.Item(2, 7).Value = 20
$dept.Cells.Item(2, 8).Value = 10

$dept.Range("A1:H2").Borders.LineStyle = 1

$dept.Hyperlinks.Add($dept.Cells.Item(2, 3), "mailto:finance@bmw.in")
$dept.Cells.Item(2, 3).Borders.LineStyle = 1

$dept.Columns.Item(1).ColumnWidth = 7
$dept.Columns.Item(2).ColumnWidth = 16
$dept.Columns.Item(3).ColumnWidth = 18.29
$dept.Columns.Item(4).ColumnWidth = 17.71
$dept.Columns.Item(5).ColumnWidth = 15
$dept.Columns.Item(6).ColumnWidth = 13.43
$dept.Columns.Item(8).ColumnWidth = 13.14

$dept.PageSetup.Orientation = 1
$dept.Range("H2").Select()

# ---------------------------------------------------------------------
# 4. Designation sheet
# ---------------------------------------------------------------------
$desig = $wb.Worksheets.Item("Designation")

$desig.Cells.Item(1, 1).Value = "Designation"
$desig.Cells.Item(1, 1).Font.Bold = $true
$desig.Cells.Item(2, 1).Value = "Quality TestEngineer"

$desig.Range("A1:A2").Borders.LineStyle = 1
$desig.Columns.Item(1).ColumnWidth = 19

$desig.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5. EmployeeGrade sheet (no bold headers, no explicit page setup)
# ---------------------------------------------------------------------
$grade = $wb.Worksheets.Item("EmployeeGrade")

$grade.Cells.Item(1, 1).Value = "Categorycode"
$grade.Cells.Item(1, 2).Value = "Categoryname"
$grade.Cells.Item(2, 1).Value = "Quad 586"
$grade.Cells.Item(2, 2).Value = "qa engineer"

$grade.Range("A1:B2").Borders.LineStyle = 1

$grade.Columns.Item(1).ColumnWidth = 12.86
$grade.Columns.Item(2).ColumnWidth = 13.71

$grade.Range("B2").Select()

# ---------------------------------------------------------------------
# 6. Passportvisaalert sheet (all cells bold, becomes the active sheet)
# ---------------------------------------------------------------------
$alert = $wb.Worksheets.Item("Passportvisaalert")

$alert.Cells.Item(1, 1).Value = " visaexpirationalert month"
$alert.Cells.Item(1, 2).Value = " Passportexpirationalert month"
$alert.Cells.Item(2, 1).Value = 10
$alert.Cells.Item(2, 2).Value = 11

$alert.Range("A1:B2").Font.Bold = $true
$alert.Range("A1:B2").Borders.LineStyle = 1

$alert.Columns.Item(1).ColumnWidth = 29.29
$alert.Columns.Item(2).ColumnWidth = 24.14

$alert.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 7. Division sheet: add a bold "Division" header row above the
#    existing list, and give the whole used range a thin border.
# ---------------------------------------------------------------------
$division = $wb.Worksheets.Item("Division")
$division.Rows.Item(1).Insert()
$division.Cells.Item(1, 1).Value = "Division"
$division.Cells.Item(1, 1).Font.Bold = $true
$division.Range("A1:A7").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 8. Passportvisaalert becomes the active/selected sheet & tab.
# ---------------------------------------------------------------------
$alert.Activate()
$alert.Range("A1").Select()
